$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the existing "Workbook" row (old row 243)
# to hold the new "WebImageCellValue" snippet entries, keeping the table
# sorted alphabetically by Class (column A).
$ws.Rows.Item(243).Insert()
$ws.Rows.Item(243).Insert()

$ws.Range("A243").Value = "WebImageCellValue"
$ws.Range("B243").Value = "address"
$ws.Range("D243").Value = "excel-data-types-web-image"
$ws.Range("E243").Value = "openImage"

$ws.Range("A244").Value = "WebImageCellValue"
$ws.Range("B244").Value = "type"
$ws.Range("D244").Value = "excel-data-types-web-image"
$ws.Range("E244").Value = "insertImage"

# Grow the Snippets table to include the two new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E304"))

# Restore the view state (frozen pane / selection) to match the edited file.
$ws.Range("A229").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E244").Select()
